$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 407, shifting existing rows 407-493 down to 408-494.
$ws.Range("A407").EntireRow.Insert()

# Fill in the new row 407 with its data. Columns A,B,C,E,F,G,H,I,N,Q,R keep
# the same repeated values as the rest of this data block.
$ws.Range("A407").Value = 3
$ws.Range("B407").Value = 'Femacal de La Calera'
$ws.Range("C407").Value = 'Coquimbo'
$ws.Range("D407").Value = 44889
$ws.Range("E407").Value = 5
$ws.Range("F407").Value = 100112017
$ws.Range("G407").Value = 'Apio'
$ws.Range("H407").Value = 'Americana (o)'
$ws.Range("I407").Value = 'Primera'
$ws.Range("J407").Value = 230
$ws.Range("K407").Value = 8000
$ws.Range("L407").Value = 8500
$ws.Range("M407").Value = 8261
$ws.Range("N407").Value = '$/docena de matas'
$ws.Range("O407").Value = 'Provincia de Limarí'
$ws.Range("P407").Value = 1377
$ws.Range("Q407").Value = 6
$ws.Range("R407").Value = 'Hortaliza'
